$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / mean calculation
$ws.Range("F2").Value = -2
$ws.Range("F6").Value = -1
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = 1
